$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right above current row 89 (shifts existing rows 89-164 down to 92-167)
$ws.Rows.Item(89).Resize(3).Insert()

# Fill in the 3 new rows with the new "Ajo" price records
# Row 89
$ws.Cells.Item(89, 1).Value = 9
$ws.Cells.Item(89, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(89, 3).Value = "Metropolitana"
$ws.Cells.Item(89, 4).Value = 44554
$ws.Cells.Item(89, 5).Value = 13
$ws.Cells.Item(89, 6).Value = 100112003
$ws.Cells.Item(89, 7).Value = "Ajo"
$ws.Cells.Item(89, 8).Value = "Rosado"
$ws.Cells.Item(89, 9).Value = "1a nueva(o)"
$ws.Cells.Item(89, 10).Value = 5000
$ws.Cells.Item(89, 11).Value = 3000
$ws.Cells.Item(89, 12).Value = 3200
$ws.Cells.Item(89, 13).Value = 3100
$ws.Cells.Item(89, 14).Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(89, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(89, 16).Value = 155
$ws.Cells.Item(89, 17).Value = 20
$ws.Cells.Item(89, 18).Value = "Hortaliza"

# Row 90
$ws.Cells.Item(90, 1).Value = 9
$ws.Cells.Item(90, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(90, 3).Value = "Metropolitana"
$ws.Cells.Item(90, 4).Value = 44554
$ws.Cells.Item(90, 5).Value = 13
$ws.Cells.Item(90, 6).Value = 100112003
$ws.Cells.Item(90, 7).Value = "Ajo"
$ws.Cells.Item(90, 8).Value = "Rosado"
$ws.Cells.Item(90, 9).Value = "2a nueva(o)"
$ws.Cells.Item(90, 10).Value = 3200
$ws.Cells.Item(90, 11).Value = 2200
$ws.Cells.Item(90, 12).Value = 2500
$ws.Cells.Item(90, 13).Value = 2350
$ws.Cells.Item(90, 14).Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(90, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(90, 16).Value = 118
$ws.Cells.Item(90, 17).Value = 20
$ws.Cells.Item(90, 18).Value = "Hortaliza"

# Row 91
$ws.Cells.Item(91, 1).Value = 9
$ws.Cells.Item(91, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(91, 3).Value = "Metropolitana"
$ws.Cells.Item(91, 4).Value = 44554
$ws.Cells.Item(91, 5).Value = 13
$ws.Cells.Item(91, 6).Value = 100112003
$ws.Cells.Item(91, 7).Value = "Ajo"
$ws.Cells.Item(91, 8).Value = "Rosado"
$ws.Cells.Item(91, 9).Value = "3a nueva (o)"
$ws.Cells.Item(91, 10).Value = 140
$ws.Cells.Item(91, 11).Value = 1500
$ws.Cells.Item(91, 12).Value = 1800
$ws.Cells.Item(91, 13).Value = 1650
$ws.Cells.Item(91, 14).Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(91, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(91, 16).Value = 82
$ws.Cells.Item(91, 17).Value = 20
$ws.Cells.Item(91, 18).Value = "Hortaliza"

# The Fecha (date) column D uses a date-formatted style (same style as the rest of column D);
# make sure the newly inserted D cells carry that style/number-format explicitly.
$ws.Range("D89:D91").NumberFormat = $ws.Range("D88").NumberFormat
